$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated team-specific transition probability values (CSUN_B matrix)
$updates = @{
    "B2" = 0.1517615176151761
    "C2" = 0.6097560975609756
    "J2" = 0.01626016260162602
    "O2" = 0.002710027100271003
    "P2" = 0.1327913279132791
    "S2" = 0.08672086720867209
    "B3" = 0.008097165991902834
    "C3" = 0.02834008097165992
    "J3" = 0.04048582995951417
    "P3" = 0.7287449392712551
    "S3" = 0.194331983805668
    "J4" = 0.06666666666666667
    "P4" = 0.5555555555555556
    "S4" = 0.3777777777777778
    "B6" = 0.06481481481481481
    "D6" = 0.009259259259259259
    "E6" = 0.00308641975308642
    "F6" = 0.07716049382716049
    "J6" = 0.2438271604938272
    "O6" = 0.0308641975308642
    "Q6" = 0.1111111111111111
    "R6" = 0.07407407407407407
    "S6" = 0.3858024691358025
    "B7" = 0.1226053639846743
    "D7" = 0.03831417624521073
    "F7" = 0.05363984674329502
    "J7" = 0.1187739463601533
    "O7" = 0.0421455938697318
    "Q7" = 0.1379310344827586
    "R7" = 0.103448275862069
    "S7" = 0.3831417624521073
    "B8" = 0.1258741258741259
    "D8" = 0.008741258741258742
    "E8" = 0.001748251748251748
    "F8" = 0.0437062937062937
    "J8" = 0.1258741258741259
    "O8" = 0.01923076923076923
    "Q8" = 0.1188811188811189
    "R8" = 0.1066433566433566
    "S8" = 0.4493006993006993
    "B9" = 0.09929078014184398
    "D9" = 0.02836879432624113
    "E9" = 0.007092198581560284
    "F9" = 0.09929078014184398
    "J9" = 0.09929078014184398
    "O9" = 0.01418439716312057
    "Q9" = 0.1276595744680851
    "R9" = 0.09929078014184398
    "S9" = 0.425531914893617
    "B10" = 0.1086519114688129
    "D10" = 0.01676727028839705
    "F10" = 0.09456740442655935
    "J10" = 0.1153588195841717
    "O10" = 0.0187793427230047
    "Q10" = 0.1891348088531187
    "R10" = 0.08517773306505701
    "S10" = 0.3715627095908786
    "G11" = 0.1895261845386534
    "J11" = 0.09226932668329177
    "K11" = 0.2418952618453865
    "L11" = 0.4663341645885287
    "S11" = 0.009975062344139651
    "G12" = 0.7074468085106383
    "J12" = 0.2393617021276596
    "K12" = 0.01063829787234043
    "S12" = 0.0425531914893617
    "G13" = 0.6588235294117647
    "J13" = 0.3058823529411765
    "S13" = 0.03529411764705882
    "F15" = 0.0375
    "H15" = 0.1625
    "I15" = 0.0375
    "J15" = 0.3416666666666667
    "K15" = 0.06666666666666667
    "M15" = 0.008333333333333333
    "O15" = 0.05
    "S15" = 0.2958333333333333
    "F16" = 0.02419354838709677
    "H16" = 0.1774193548387097
    "I16" = 0.0846774193548387
    "J16" = 0.3951612903225806
    "K16" = 0.1008064516129032
    "M16" = 0.04435483870967742
    "O16" = 0.03225806451612903
    "S16" = 0.1411290322580645
    "F17" = 0.02237136465324385
    "H17" = 0.174496644295302
    "I17" = 0.04697986577181208
    "J17" = 0.4340044742729307
    "K17" = 0.09172259507829977
    "M17" = 0.029082774049217
    "O17" = 0.06935123042505593
    "S17" = 0.1319910514541387
    "F18" = 0.02745098039215686
    "H18" = 0.2156862745098039
    "I18" = 0.05490196078431372
    "J18" = 0.4117647058823529
    "K18" = 0.1137254901960784
    "M18" = 0.0196078431372549
    "O18" = 0.04313725490196078
    "S18" = 0.1137254901960784
    "F19" = 0.02530044275774826
    "H19" = 0.228336495888678
    "I19" = 0.04743833017077799
    "J19" = 0.3523086654016445
    "K19" = 0.1214421252371917
    "M19" = 0.03542061986084757
    "N19" = 0.001265022137887413
    "O19" = 0.058191018342821
    "S19" = 0.1302972802024036
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
